$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORT_DATE (text, not a real date serial - matches existing inline-string format)
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# Numeric financial figures that changed
$ws.Range("O2").Value = 393797078.78
$ws.Range("P2").Value = 121527467.58
$ws.Range("Q2").Value = 10288960.37
$ws.Range("S2").Value = 56985757.13
$ws.Range("U2").Value = 94749305.31999999
$ws.Range("W2").Value = 173026372.18
$ws.Range("X2").Value = 91015223.28
$ws.Range("Z2").Value = 8689836.07
$ws.Range("AB2").Value = 220770706.6
$ws.Range("AF2").Value = 121.6439117345
$ws.Range("AG2").Value = 43.9379521849

# Ratio columns cleared out (now blank cells) in the new data
$ws.Range("R2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
